$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("backup_2021_2")

# Row 2
$ws.Range("B2").Value = ""
$ws.Range("F2").Value = "Kovács Gusztáv Márk"
$ws.Range("G2").Value = "Ördög Márk"
$ws.Range("H2").Value = "Siklósi Balázs"

# Row 3
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = ""

# Row 7
$ws.Range("G7").Value = "Hajdu Krisztián"

# Row 8
$ws.Range("D8").Value = ""
